# Enhanced py files and add new col to files tab
#
# Updates the four SQL "TabQuery"/"StatQuery" cells on Sheet1 with the
# enhanced queries (new joins against df_diagnosis / df_genomic_info,
# phs004231 -> phs002431, extra "Library Strategy" / "File Type" columns
# on the Files tab query, LIMIT clauses, etc.), keeps the font/wrap
# formatting on those cells, widens row heights where the new text needs
# more room, and nudges the visible top row of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ParticipantsTab / TabQuery (B2): add diagnosis/program/file/genomic joins ---
$participantsTabQuery = @'
WITH Distinct_Samples AS (
    SELECT DISTINCT
        sp.participant_id,
        sp.study_participant_id,
        s.study_name,
        s.phs_accession,
        sp.gender,
        smp.sample_id
    FROM 
        df_participant sp
    JOIN 
        df_study s ON sp."study.phs_accession" = s.phs_accession
    JOIN 
        df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
    JOIN
        df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
    JOIN
        df_program p ON p.program_acronym = s."program.program_acronym"
    JOIN
        df_file f1 ON f1."sample.sample_id" = smp.sample_id
    JOIN
        df_genomic_info gi ON gi."file.file_id" = f1.file_id
    WHERE 
        s.phs_accession = 'phs002431' AND sp.gender = 'Female'
),
Sample_Limit AS (
    SELECT
        participant_id,
        study_participant_id,
        study_name,
        phs_accession,
        gender,
        sample_id,
        ROW_NUMBER() OVER (PARTITION BY study_participant_id ORDER BY sample_id) as rn
    FROM 
        Distinct_Samples
)
SELECT
    participant_id AS "Participant ID", 
    study_name AS "Study Name",
    phs_accession AS Accession,
    gender AS Gender,
    GROUP_CONCAT(
        CASE 
            WHEN rn <= 5 THEN sample_id 
        END, ', '
    ) ||
    CASE 
        WHEN MAX(rn) > 5 THEN ', ...' 
        ELSE '' 
    END AS Samples
FROM 
    Sample_Limit
GROUP BY
    participant_id, 
    study_name,
    phs_accession,
    gender
LIMIT 100;
'@

# --- ParticipantsTab / StatQuery (C2): add diagnosis/genomic joins, fix accession ---
$participantsStatQuery = @'
SELECT
    COUNT(DISTINCT s.study_name) AS "Studies",
    COUNT(DISTINCT sp.participant_id) AS "Participants",
    COUNT(DISTINCT smp.sample_id) AS "Samples",
    COUNT(DISTINCT f.file_id) AS "Files"
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN 
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN 
    df_file f ON f."sample.sample_id" = smp.sample_id
JOIN 
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN 
    df_genomic_info gi ON gi."file.file_id" = f.file_id
WHERE 
    s.phs_accession = 'phs002431' AND sp.gender = 'Female';
'@

# --- SamplesTab / TabQuery (B3): add diagnosis/program/file/genomic joins ---
$samplesTabQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession,
    smp.sample_tumor_status AS Tumor,
    smp.sample_type AS "Analyte Type"
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
    s.phs_accession = 'phs002431' AND sp.gender = 'Female'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

# --- FilesTab / TabQuery (B4): new "Library Strategy" column + diagnosis join ---
$filesTabQuery = @'
SELECT
    f1.file_name AS "File Name",
    s.study_name AS "Study Name",
    s.phs_accession AS "Accession",
    sp.participant_id AS "Participant Id",
    COALESCE((
        SELECT
            REPLACE(GROUP_CONCAT(CASE WHEN rn <= 5 THEN smp.sample_id ELSE NULL END, ', '), ', , ', ', ') ||
            CASE WHEN MAX(rn) > 5 THEN ', ...' ELSE '' END
        FROM (
            SELECT
                smp.sample_id,
                ROW_NUMBER() OVER (ORDER BY smp.sample_id) AS rn
            FROM df_sample smp
            WHERE smp."participant.study_participant_id" = sp.study_participant_id
        ) smp
    ), '') AS "Sample Id",
    f1.file_type AS "File Type",
    gi.library_strategy AS "Library Strategy"
FROM 
    df_study s
INNER JOIN 
    df_participant sp ON sp."study.phs_accession" = s.phs_accession
INNER JOIN  
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
INNER JOIN 
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
INNER JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
INNER JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
INNER JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
WHERE 
    s.phs_accession = 'phs002431' AND sp.gender = 'Female'
GROUP BY
    f1.file_name,
    s.study_name,
    s.phs_accession,
    sp.participant_id,
    f1.file_type,
    gi.library_strategy
ORDER BY 
    f1.file_name ASC
LIMIT 100;
'@

# Write the enhanced query text back into the four query cells.
$ws.Range("B2").Value = $participantsTabQuery
$ws.Range("C2").Value = $participantsStatQuery
$ws.Range("B3").Value = $samplesTabQuery
$ws.Range("B4").Value = $filesTabQuery

# Re-assert the existing formatting (12pt Calibri, wrapped text) on the
# updated cells so the longer text still displays correctly.
$queryCells = $ws.Range("B2,C2,B3,B4")
$queryCells.Font.Size = 12
$queryCells.WrapText = $true

# The longer queries need taller rows to show fully (row 3 grows to the
# worksheet's max row height; rows 2 and 4 already were at the max).
$ws.Rows(2).RowHeight = 409.6
$ws.Rows(3).RowHeight = 409.6
$ws.Rows(4).RowHeight = 409.6

# Scroll the sheet's visible top row down one row.
$excel.ActiveWindow.ScrollRow = 3
